# B6-PowerPoint.pptx edit
#
# 1) Three tables (on the slides that show the balance-sheet / accounts
#    exercises) get their table style switched from the custom
#    "Table_0" style ({7691C40B-05DE-45CB-BEDE-2EF4C9D3A5B9}) to the
#    built-in style {CBA3ED6C-D994-4727-B0B4-6AE3C4B793EA}.
#
# 2) The deck's theme is swapped from the "Integral / Red Violet" theme
#    back to the default "Office" colour scheme (dk1/lt1/dk2/lt2/accentN/
#    hlink/folHlink) - i.e. the Design the deck uses switches from the
#    pink/purple "Red Violet" palette to the plain blue/orange "Office"
#    palette.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-apply the table style on the three slides that contain a table.
# ---------------------------------------------------------------------
$newTableStyleId = "{CBA3ED6C-D994-4727-B0B4-6AE3C4B793EA}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Swap the theme colours back to the stock "Office" palette.
# ---------------------------------------------------------------------
# COM ThemeColor.RGB uses the Windows COLORREF byte order (0x00BBGGRR),
# so build each value from the target RRGGBB hex by reversing the bytes.
function ColorRefFromRGBHex([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Slot order exposed by ThemeColorScheme.Colors(1..12):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1 .. 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$themeColorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColorScheme.Count; $i++) {
    $themeColorScheme.Colors($i).RGB = ColorRefFromRGBHex($officeColors[$i - 1])
}

$p.Save()
